$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bom")

# "fixed eeprom in bom": the RS-Components part number for the EEPROM
# (CAT24C256, row 13) was wrong. Correct it from 800-9812 to 808-0199.
$ws.Range("C13").Value = "808-0199"

# The cell's hyperlink needs to keep pointing at the (unchanged) RS-Components
# product page, but must no longer show the stale "800-9812" custom display
# text - Excel drops the custom display text once it no longer matches.
# This runtime only supports adding/removing whole hyperlinks (in-place
# property edits on a fetched Hyperlink object create a duplicate instead of
# mutating it), so rebuild the sheet's hyperlinks from scratch: drop them all
# and re-add each with its original target, restoring the original display
# text everywhere except C13.
$ws.Range("C5").Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("C2"),  "https://at.rs-online.com/web/p/leiterplatten-buchsen/6742369/", "", "", " 674-2369") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"),  "https://at.rs-online.com/web/p/leiterplatten-header/6812979", "", "", " 681-2979") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"),  "https://at.rs-online.com/web/p/leiterplatten-header/6812975", "", "", " 681-2975") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"),  "https://at.rs-online.com/web/p/smd-widerstande/8075380/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"),  "https://at.rs-online.com/web/p/smd-widerstande/8075550/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"),  "https://at.rs-online.com/web/p/smd-widerstande/8075579/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"),  "https://at.rs-online.com/web/p/keramik-vielschichtkondensatoren/1721482", "", "", " 172-1482") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"),  "https://at.rs-online.com/web/p/led/8610100/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), "https://at.rs-online.com/web/p/tastschalter/1359508", "", "", " 135-9508") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "https://at.rs-online.com/web/p/temperatursensoren-und-feuchtigkeitssensoren/5402849", "", "", " 540-2849") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "https://at.rs-online.com/web/p/a-d-wandler/6696058", "", "", " 669-6058") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "https://at.rs-online.com/web/p/eeproms/8009812") | Out-Null

# Adding a hyperlink with a TextToDisplay overwrites the cell's literal text
# with that display string (here with its leading space) - restore every
# cell's real (space-free) text content afterwards.
$ws.Range("C2").Value  = "674-2369"
$ws.Range("C3").Value  = "681-2979"
$ws.Range("C4").Value  = "681-2975"
$ws.Range("C5").Value  = "807-5380"
$ws.Range("C6").Value  = "807-5550"
$ws.Range("C7").Value  = "807-5579"
$ws.Range("C8").Value  = "172-1482"
$ws.Range("C9").Value  = "861-0100"
$ws.Range("C10").Value = "135-9508"
$ws.Range("C11").Value = "540-2849"
$ws.Range("C12").Value = "669-6058"
$ws.Range("C13").Value = "808-0199"

# Reflect the last-edited cell in the saved selection.
$ws.Range("C13").Select()
